$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 10073.5
$ws.Range("I4").Value = 10073.5
$ws.Range("K4").Value = 10073.5
$ws.Range("M4").Value = -9959.5

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H15").Value = 360
$ws.Range("I15").Value = 360
$ws.Range("K15").Value = 1080
$ws.Range("M15").Value = -911

$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 30000
$ws.Range("K34").Value = 30000
$ws.Range("M34").Value = -29797

$ws.Range("H36").Value = 30000
$ws.Range("I36").Value = 30000
$ws.Range("K36").Value = 30000
$ws.Range("M36").Value = -29285

$ws.Range("H43").Value = 2666.6667
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H62").Value = 4385.7144
$ws.Range("I62").Value = 4116.6665
$ws.Range("K62").Value = 4116.6665
$ws.Range("M62").Value = -3492.6665

$ws.Range("H65").Value = 4385.7144
$ws.Range("I65").Value = 4116.6665
$ws.Range("K65").Value = 20583.3325
$ws.Range("M65").Value = -17463.3325

$ws.Range("H101").Value = 1831.6666
$ws.Range("I101").Value = 1497.5
$ws.Range("J101").Value = 2500
$ws.Range("K101").Value = 4492.5
$ws.Range("L101").Value = 7500
$ws.Range("M101").Value = -2870.5
$ws.Range("N101").Value = -10744

$ws.Range("H112").Value = 5124.6
$ws.Range("J112").Value = 5124.6
$ws.Range("L112").Value = 15373.8
$ws.Range("N112").Value = -17589.8

$ws.Range("H116").Value = 4475
$ws.Range("J116").Value = 3950
$ws.Range("L116").Value = 3950
$ws.Range("N116").Value = -10834

$ws.Range("H138").Value = 6518.9
$ws.Range("I138").Value = 7344.2
$ws.Range("J138").Value = 5693.6
$ws.Range("K138").Value = 22032.6
$ws.Range("L138").Value = 17080.8
$ws.Range("M138").Value = -16892.6
$ws.Range("N138").Value = -27360.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 47.333332
$ws.Range("I5").Value = 47
$ws.Range("K5").Value = 47
$ws.Range("M5").Value = 65

$ws.Range("H28").Value = 2381.8572
$ws.Range("I28").Value = 1612.1666
$ws.Range("K28").Value = 1612.1666
$ws.Range("M28").Value = -1420.1666

$ws.Range("H32").Value = 14392.267
$ws.Range("I32").Value = 11690.4
$ws.Range("J32").Value = 19796
$ws.Range("K32").Value = 11690.4
$ws.Range("L32").Value = 19796
$ws.Range("M32").Value = -11403.4
$ws.Range("N32").Value = -20370

$ws.Range("H99").Value = 2381.8572
$ws.Range("I99").Value = 1612.1666
$ws.Range("K99").Value = 1612.1666
$ws.Range("M99").Value = 1382.8334

$ws.Range("H110").Value = 1149.8572
$ws.Range("I110").Value = 1149.8572
$ws.Range("K110").Value = 1149.8572
$ws.Range("M110").Value = 895.1428000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 47.333332
$ws.Range("I4").Value = 47
$ws.Range("K4").Value = 47
$ws.Range("M4").Value = 68

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H75").Value = 3382.5
$ws.Range("I75").Value = 3382.5
$ws.Range("K75").Value = 3382.5
$ws.Range("M75").Value = -2446.5

$ws.Range("H78").Value = 3382.5
$ws.Range("I78").Value = 3382.5
$ws.Range("K78").Value = 10147.5
$ws.Range("M78").Value = -5467.5

$ws.Range("H82").Value = 6196.6
$ws.Range("I82").Value = 6196.6
$ws.Range("K82").Value = 6196.6
$ws.Range("M82").Value = -5813.6

$ws.Range("H85").Value = 6196.6
$ws.Range("I85").Value = 6196.6
$ws.Range("K85").Value = 6196.6
$ws.Range("M85").Value = -4870.6

$ws.Range("H97").Value = 23449.5
$ws.Range("I97").Value = 23449.5
$ws.Range("K97").Value = 23449.5
$ws.Range("M97").Value = -22458.5

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 125
$ws.Range("I22").Value = 50
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = -900

$ws.Range("H31").Value = 18771.875
$ws.Range("I31").Value = 10827.385
$ws.Range("K31").Value = 10827.385
$ws.Range("M31").Value = -10532.385

$ws.Range("H34").Value = 18771.875
$ws.Range("I34").Value = 10827.385
$ws.Range("K34").Value = 10827.385
$ws.Range("M34").Value = -10625.385

$ws.Range("H99").Value = 2510.75
$ws.Range("I99").Value = 2510.75
$ws.Range("K99").Value = 2510.75
$ws.Range("M99").Value = -1012.75

$ws.Range("H126").Value = 2510.75
$ws.Range("I126").Value = 2510.75
$ws.Range("K126").Value = 7532.25
$ws.Range("M126").Value = -5062.25

$ws.Range("H133").Value = 50765
$ws.Range("I133").Value = 15296
$ws.Range("K133").Value = 15296
$ws.Range("M133").Value = -12766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2896.077
$ws.Range("I131").Value = 2193.75
$ws.Range("K131").Value = 6581.25
$ws.Range("M131").Value = -1541.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 850
$ws.Range("J27").Value = 850
$ws.Range("L27").Value = 850
$ws.Range("N27").Value = -1182

$ws.Range("H97").Value = 2109.5
$ws.Range("I97").Value = 2108
$ws.Range("K97").Value = 2108
$ws.Range("M97").Value = -1612

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 3355.6667
$ws.Range("I132").Value = 2033.5
$ws.Range("K132").Value = 6100.5
$ws.Range("M132").Value = -3570.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800

$ws.Range("H27").Value = 800

$ws.Range("H55").Value = 7133.1665
$ws.Range("I55").Value = 6949.5
$ws.Range("K55").Value = 6949.5
$ws.Range("M55").Value = -6776.5

$ws.Range("H61").Value = 2969.3333
$ws.Range("I61").Value = 954
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 954
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -752
$ws.Range("N61").Value = -7404

$ws.Range("H68").Value = 2638.2307
$ws.Range("I68").Value = 2921.889
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2921.889
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -2172.889
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 2638.2307
$ws.Range("I71").Value = 2921.889
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 14609.445
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -10865.445
$ws.Range("N71").Value = -17488

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H100").Value = 4300
$ws.Range("I100").Value = 3816.6667
$ws.Range("K100").Value = 3816.6667
$ws.Range("M100").Value = -3275.6667

$ws.Range("H113").Value = 2969.3333
$ws.Range("I113").Value = 954
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 954
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 1216
$ws.Range("N113").Value = -11340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 24000
$ws.Range("I61").Value = 24000
$ws.Range("K61").Value = 24000
$ws.Range("M61").Value = -23708

$ws.Range("H62").Value = 2743.5
$ws.Range("I62").Value = 2743.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2743.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2119.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2743.5
$ws.Range("I65").Value = 2743.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13717.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10597.5
$ws.Range("N65").ClearContents()

$ws.Range("H131").Value = 42000
$ws.Range("J131").Value = 42000
$ws.Range("L131").Value = 42000
$ws.Range("N131").Value = -52080
